# DAY 08 - excel automation day 2
#
# Renames the worksheet ("data" -> "Sheet1") and moves the active
# selection from A3 to B5, matching the saved-workbook-state changes
# captured in the commit's OOXML diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet tab.
$ws.Name = "Sheet1"

# Move the selection to B5 (becomes the new active cell / sqref).
$ws.Range("B5").Select()
